$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the lesson grid (A1:F11) ---
# Row 1: section headers
$ws.Range("A1").Value = "Numbers"
$ws.Range("B1").Value = "Math"
$ws.Range("C1").Value = "Letters"
$ws.Range("D1").Value = "Words"
$ws.Range("E1").Value = "Shapes"
$ws.Range("F1").Value = "Colors"

# Row 2
$ws.Range("A2").Value = "0, zero, ze-ro"
$ws.Range("B2").Value = "Addition,1 + 1 = 2,1 + 2 = 3,1 + 3 = 4"
$ws.Range("C2").Value = "A,Ant"

# Row 3
$ws.Range("A3").Value = "1, one, wun"
$ws.Range("B3").Value = "Commutative Property,2 + 1 = 1 + 2 = 3"
$ws.Range("C3").Value = "B,Bird"

# Row 4
$ws.Range("A4").Value = "2, two, tu"
$ws.Range("B4").Value = "Identity Property,3 + 0 = 3,0 + 5 = 5"
$ws.Range("C4").Value = "C,Cat"

# Row 5
$ws.Range("A5").Value = "3, three, three"
$ws.Range("B5").Value = "Associative Property,(2 + 1) + 3 = 2 + (1 + 3) = 6"
$ws.Range("C5").Value = "D,Dog"

# Row 6
$ws.Range("A6").Value = "4, four, for"
$ws.Range("B6").Value = "Examples,4 + 4 = ?, 5 + 1 = ?,6 + 0 = ?"
$ws.Range("C6").Value = "E,Elephant"

# Row 7
$ws.Range("A7").Value = "5, five, fayv"
$ws.Range("B7").Value = "Subtraction,3 - 2 = 1,5 - 2 = 3,2 -2 = 0"
$ws.Range("C7").Value = "F,Frog"

# Row 8
$ws.Range("A8").Value = "6, six, six"
$ws.Range("B8").Value = "Subtraction,3 - 2 = 1,3 is the minuend,2 is the subtrahend,1 is the difference"
$ws.Range("C8").Value = "G,Giraffe"

# Row 9
$ws.Range("A9").Value = "7, seven, se-ven"
$ws.Range("B9").Value = "Identity Property,5 - 0 = 5, 9 - 0 = 9"
$ws.Range("C9").Value = "H,Horse"

# Row 10
$ws.Range("A10").Value = "8, eight, eyt"
$ws.Range("B10").Value = "Commutative Property,2 + 1 = 1 + 2 = 3"
$ws.Range("C10").Value = "I,Iguana"

# Row 11
$ws.Range("A11").Value = "9, nine, nayn"
$ws.Range("B11").Value = "Examples,3 - 2 = ?,5 - 4 = ?,6 - 0 = ?"
$ws.Range("C11").Value = "J,Jellyfish"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 16.453125
$ws.Columns.Item(2).ColumnWidth = 45.54296875

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 19
$ws.Range("A2:A11").EntireRow.RowHeight = 64.5

# --- Alignment / wrap for the whole used range ---
$full = $ws.Range("A1:F11")
$full.HorizontalAlignment = -4131
$full.VerticalAlignment = -4160
$full.WrapText = $true

# --- Selection ---
$ws.Range("B5").Select()
